# Insert a new student row ("Дзівідзінська Мар'яна") above "Орлов Володимир",
# pushing the remaining rows down by one, and renumber the A column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing values first (Value2 reads reliably; Value is used for writes).
$row3B = $ws.Cells.Item(3, 2).Value2
$row3C = $ws.Cells.Item(3, 3).Value2
$row4B = $ws.Cells.Item(4, 2).Value2
$row4C = $ws.Cells.Item(4, 3).Value2

# Shift the last existing data row (old row 4: Ямковий Андрій) down to row 5.
# Copy the numbering-column format (bold/border/centered) from A4 onto A5.
$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = $row4B
$ws.Cells.Item(5, 3).Value = $row4C

# Shift the old row 3 (Орлов Володимир) down to row 4.
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = $row3B
$ws.Cells.Item(4, 3).Value = $row3C

# Write the new student into row 3.
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Дзівідзінська Мар'яна"
$ws.Cells.Item(3, 3).Value = 0
